# The image filenames (column L) and their associated stats (columns M:V)
# were reassigned to different rows ("cleaned up the image folder" / image
# renames). Apply the row permutation by rotating the L:V block of values
# around the cycle of affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycle order: row at position i receives the L:V content that currently
# lives in the row at position i+1 (wrapping around).
$rows = @(7, 16, 21, 40, 20, 19, 23, 15, 12, 18, 28)

# Snapshot current L:V values (as a list of column values) for every row
# in the cycle before any writes happen. (Use Value2 -- Value's getter is
# unreliable in this host for reading back ranges.)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = $ws.Range("L$r`:V$r").Value2
}

# Write each row's new L:V block = snapshot of the next row in the cycle.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $destRow = $rows[$i]
    $srcRow = $rows[($i + 1) % $rows.Length]
    $ws.Range("L$destRow`:V$destRow").Value2 = $snapshot[$srcRow]
}
